# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# Mirrors the commit: "Added team record to data" - W/L/T live on the same
# sheet (columns AD:AF) rather than a separate sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 44

# --- Header row (row 1): new column headers in AD1:AF1 -----------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the existing header row (bold, centered/top aligned,
# thin box border) for the three new header cells.
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

# --- Data rows (2..44): every player row gets the team's record --------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 88   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
